# Swap the contents of column C (codeforiati:group-code) and
# column D (codeforiati:group-name), including the header row.
# This reproduces the shared-string reordering seen in the diff,
# where each (code, name) pair of <si> entries was swapped so that
# the "name" string now precedes the "code" string in the table,
# while every cell that used to show the code now shows the name
# and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
